# Update to new files
# - Add "U" (unknown sex marker) into column D (sex) for each data row.
# - Clear the sire/dam tallies previously held in columns F and G.
# - Move the active selection to D7 (just below the last data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column D (sex) with "U" for each data row (2-6)
$ws.Range("D2:D6").Value = "U"

# Clear out the old sire/dam columns (F and G) for the data rows
$ws.Range("F2:G6").ClearContents()

# Update the saved selection to D7
$ws.Range("D7").Select()
